$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 4: "[OARC 31] ... Austin, TX, USA, October 2019." -- the stray
# "_GoBack" bookmark that used to sit between "...October 2019" and the
# trailing "." is being relocated (see hunk 1 below), so here we just
# fold the trailing period into the same run as the rest of the date,
# which drops the now-unwanted bookmark markers in the process.
# ---------------------------------------------------------------------
$oldAustin = " Austin, TX, USA, October 2019."
$newAustin = " Austin, TX, USA, October 2019."
$foundAustin = $d.Content.Find.Execute($oldAustin, $false, $false, $false, $false, $false, $true, 1, $false, $newAustin, 2)
Write-Host "Austin fix found: $foundAustin"

# ---------------------------------------------------------------------
# Hunk 1: "[NDSS '21] ... To appear in the Network and Distributed
# System Security Symposium 2021, Virtual event, February 2021." ->
# "... In Proceedings of the Network and Distributed System Security
# Symposium 2021, Virtual event, February 2021." with a "_GoBack"
# bookmark inserted right after "In Proceedings of".
# ---------------------------------------------------------------------
$oldNdss = "To appear in the Network and Distributed System Security Symposium 2021, Virtual event, February 2021."
$newNdss = "In Proceedings of the Network and Distributed System Security Symposium 2021, Virtual event, February 2021."
$foundNdss = $d.Content.Find.Execute($oldNdss, $false, $false, $false, $false, $false, $true, 1, $false, $newNdss, 2)
Write-Host "NDSS text found: $foundNdss"

$afterInProc = $d.Content.Duplicate
$foundInProc = $afterInProc.Find.Execute("In Proceedings of", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "'In Proceedings of' found: $foundInProc"
$bmRange = $d.Range($afterInProc.End, $afterInProc.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# Hunk 2: "An End-to-End, ... How Far Have We [gramStart]Come?[/gramStart]"
# -> merge into a single run "An End-to-End, ... How Far Have We Come?"
# and drop the proofErr markers around "Come?" and the following ".".
# ---------------------------------------------------------------------
$oldDoe = "An End-to-End, Large-Scale Measurement of DNS-over-Encryption: How Far Have We Come?"
$newDoe = "An End-to-End, Large-Scale Measurement of DNS-over-Encryption: How Far Have We Come?"
$foundDoe = $d.Content.Find.Execute($oldDoe, $false, $false, $false, $false, $false, $true, 1, $false, $newDoe, 2)
Write-Host "DoE title found: $foundDoe"

# The trailing "[gramEnd]" marker sits between the "." run and the bold
# " " run that follows it; touch that span (restoring the original
# bold/non-bold split afterwards) so the marker gets dropped too.
$doeTail = $d.Content.Duplicate
$foundDoeTail = $doeTail.Find.Execute("Come?. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "DoE tail found: $foundDoeTail"
$dotSpace = $d.Range($doeTail.End - 2, $doeTail.End)
$dotSpace.Text = ".X"
$spaceFix = $d.Range($dotSpace.End - 1, $dotSpace.End)
$spaceFix.Text = " "
$spaceFix2 = $d.Range($dotSpace.End - 1, $dotSpace.End)
$spaceFix2.Font.Bold = 1

# ---------------------------------------------------------------------
# Hunk 3: "A Reexamination of Internationalized Domain Names:
# [gramStart]the[/gramStart] Good, the Bad and the Ugly" -> merge into a
# single run and drop the proofErr markers around "the".
# ---------------------------------------------------------------------
$oldIdn = "A Reexamination of Internationalized Domain Names: the Good, the Bad and the Ugly"
$newIdn = "A Reexamination of Internationalized Domain Names: the Good, the Bad and the Ugly"
$foundIdn = $d.Content.Find.Execute($oldIdn, $false, $false, $false, $false, $false, $true, 1, $false, $newIdn, 2)
Write-Host "IDN title found: $foundIdn"
